$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Under "Requisitos:" (A24), the three requirement rows (25-27) get re-ordered:
# the "LOM3246 - Técnicas de Caracterização de Materiais (Indicação de Conjunto)"
# entry now comes first, pushing the other two requirements down by one row
# (their own relative order stays the same).
$ws.Range("B25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$ws.Range("C25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"

$ws.Range("B26").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("C26").Value = "LOB1021 -  Física IV  (Requisito)`n"

$ws.Range("B27").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C27").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
